$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(4, 7).Value = "Dr. Shimaa Ashraf, Dr. Hend Farid, Dr. Amal Awwad, Dr. Aya Saeed, Dr. Mariam Nour El-Din"
$ws.Cells.Item(5, 7).Value = "D Wessam Atef, Dr. Amal Awwad, Dr. Sara Nabil, Dr. Nourhan Mohammad"
$ws.Cells.Item(6, 7).Value = "Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Mai Mustafa, Dr. Basma Hamed, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed"
$ws.Cells.Item(8, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef"
$ws.Cells.Item(9, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat"
$ws.Cells.Item(10, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(14, 7).Value = "Dr. Marian Samir, Dr. Nourhan Mohammad"
$ws.Cells.Item(15, 7).Value = "Dr. Afaf Abdallah, Dr. Marian Samir, Dr. Nourham Mostafa, Dr. Nourhan Mohammad, Dr. Ahmad Mostafa"
$ws.Cells.Item(16, 7).Value = "Dr. Manarst Al-Eslam, Dr. Walaa Ghanima, Dr. Rada Rabea, Dr. Marian Samir"
$ws.Cells.Item(17, 7).Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Youstina Magdy, Dr. Nardine, Dr. Monica"
$ws.Cells.Item(18, 7).Value = "Dr. Maryam Ashraf, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Monica, Dr. Remon, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Aya Emad"
$ws.Cells.Item(19, 7).Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(21, 7).Value = "Dr. Shimaa Ashraf, Dr. Hend Farid, Dr. Amal Awwad, Dr. Aya Saeed, Dr. Mariam Nour El-Din"
$ws.Cells.Item(22, 7).Value = "D Wessam Atef, Dr. Amal Awwad, Dr. Sara Nabil, Dr. Nourhan Mohammad"
$ws.Cells.Item(23, 7).Value = "Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Mai Mustafa, Dr. Basma Hamed, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed"
$ws.Cells.Item(25, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef"
$ws.Cells.Item(26, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat"
$ws.Cells.Item(27, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(30, 7).Value = "Dr. Mariam Gamal Sanad, Dr. Sarah Mahdy"
$ws.Cells.Item(31, 7).Value = "Dr. Marian Samir, Dr. Nourhan Mohammad"
$ws.Cells.Item(32, 7).Value = "Dr. Afaf Abdallah, Dr. Marian Samir, Dr. Nourham Mostafa, Dr. Nourhan Mohammad, Dr. Ahmad Mostafa"
$ws.Cells.Item(33, 7).Value = "Dr. Manarst Al-Eslam, Dr. Walaa Ghanima, Dr. Rada Rabea, Dr. Marian Samir"
$ws.Cells.Item(34, 7).Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Youstina Magdy, Dr. Nardine, Dr. Monica"
$ws.Cells.Item(35, 7).Value = "Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Aya Emad"
$ws.Cells.Item(36, 7).Value = "Dr. Eman Tantawi, Administrator, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(37, 7).Value = "Administrator, Dr. Kerelos Zareef, Dr. Nada Mohammad"
$ws.Cells.Item(39, 7).Value = "Dr. Omnia Mohammad, Dr. Shimaa Ashraf"
$ws.Cells.Item(40, 7).Value = "Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Mai Mustafa, Dr. Basma Hamed, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed"
$ws.Cells.Item(43, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat"
$ws.Cells.Item(44, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(45, 7).Value = "Dr. Rania Ahmad Youssef, Administrator, Dr. Mohammad Safwat"
$ws.Cells.Item(48, 7).Value = "Dr. Aya Alaa-Eldein, Dr. Afaf Abdallah, Dr. Marian Samir"
$ws.Cells.Item(50, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Manarst Al-Eslam, Dr. Aya Alaa-Eldein"
$ws.Cells.Item(51, 7).Value = "Dr. Yasmin, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Monica, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(52, 7).Value = "Dr. Yasmin, Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Cells.Item(53, 7).Value = "Dr. Eman Tantawi, Administrator, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(54, 7).Value = "Administrator, Dr. Kerelos Zareef, Dr. Nada Mohammad"
$ws.Cells.Item(56, 7).Value = "Dr. Omnia Mohammad, Dr. Shimaa Ashraf"
$ws.Cells.Item(57, 7).Value = "Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Mai Mustafa, Dr. Basma Hamed, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed"
$ws.Cells.Item(60, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat"
$ws.Cells.Item(61, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(62, 7).Value = "Dr. Rania Ahmad Youssef, Administrator, Dr. Mohammad Safwat"
$ws.Cells.Item(65, 7).Value = "Dr. Aya Alaa-Eldein, Dr. Afaf Abdallah, Dr. Marian Samir"
$ws.Cells.Item(67, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Manarst Al-Eslam, Dr. Aya Alaa-Eldein"
$ws.Cells.Item(68, 7).Value = "Dr. Yasmin, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Monica, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(69, 7).Value = "Dr. Yasmin, Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Cells.Item(70, 7).Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(72, 7).Value = "D Wessam Atef, Dr. Shimaa Ashraf, Dr. Safa Hany, Dr. Omnia Mohammad, Dr. Mariam Nour El-Din"
$ws.Cells.Item(73, 7).Value = "Dr. Shimaa Ashraf, Dr. Hend Farid, Dr. Amal Awwad, Dr. Aya Saeed, Dr. Mariam Nour El-Din"
$ws.Cells.Item(74, 7).Value = "D Wessam Atef, Dr. Sara Nabil, Dr. Omnia Mohammad, Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Cells.Item(75, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Eman M. Elsaid, Dr. Merna Said, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed"
$ws.Cells.Item(76, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mariam Toma Gerges, Dr. Mohammad Safwat"
$ws.Cells.Item(77, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mariam Toma Gerges, Dr. Mohammad Safwat"
$ws.Cells.Item(78, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat, Dr. Mayar Ahmad Embaby"
$ws.Cells.Item(79, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(83, 7).Value = "Dr. Youstina Ibrahim, Dr. Afaf Abdallah, Dr. Marian Samir"
$ws.Cells.Item(84, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Marian Samir, Dr. Manarst Al-Eslam, Dr. Aya Alaa-Eldein"
$ws.Cells.Item(85, 7).Value = "Dr. Maryam Ashraf, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Nahla, Dr. Neveen Nashaat, Dr. Monica, Dr. Wafaa Ebida, Dr. Aya Emad"
$ws.Cells.Item(86, 7).Value = "Dr. Maryam Ashraf, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Monica, Dr. Remon, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Aya Emad"
$ws.Cells.Item(87, 7).Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(88, 7).Value = "Dr. Nada Mohammad, Dr. Fatma Elhady"
$ws.Cells.Item(89, 7).Value = "Dr. Shimaa Ashraf, Dr. Hend Farid, Dr. Amal Awwad, Dr. Aya Saeed, Dr. Mariam Nour El-Din"
$ws.Cells.Item(90, 7).Value = "D Wessam Atef, Dr. Sara Nabil, Dr. Omnia Mohammad, Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Cells.Item(91, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Eman M. Elsaid, Dr. Merna Said, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed"
$ws.Cells.Item(92, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mariam Toma Gerges, Dr. Mohammad Safwat"
$ws.Cells.Item(93, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat"
$ws.Cells.Item(94, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(98, 7).Value = "Dr. Afaf Abdallah, Dr. Walaa Ghanima, Dr. Nourhan Hosni, Dr. Nourhan Mohammad"
$ws.Cells.Item(101, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Marian Samir, Dr. Manarst Al-Eslam, Dr. Aya Alaa-Eldein"
$ws.Cells.Item(102, 7).Value = "Dr. Maryam Ashraf, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Nahla, Dr. Neveen Nashaat, Dr. Monica, Dr. Wafaa Ebida, Dr. Aya Emad"
$ws.Cells.Item(103, 7).Value = "Dr. Maryam Ashraf, Dr. Ola Abd Al-Fattah, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Gehad Salah, Dr. Yassmen Ahmad, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Aya Hanafy"
$ws.Cells.Item(104, 7).Value = "Dr. Eman Tantawi, Administrator, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(106, 7).Value = "D Wessam Atef, Dr. Amal Awwad, Dr. Sara Nabil, Dr. Nourhan Mohammad"
$ws.Cells.Item(107, 7).Value = "Dr. Amal Awwad, Dr. Sara Nabil, Dr. Nourhan Mohammad"
$ws.Cells.Item(108, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Basma Hamed, Dr. Arwa Al-Sayed"
$ws.Cells.Item(111, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat"
$ws.Cells.Item(112, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(115, 7).Value = "Dr. Mariam Gamal Sanad, Dr. Sarah Mahdy"
$ws.Cells.Item(116, 7).Value = "Dr. Enas Omran, Dr. Afaf Abdallah, Dr. Nourham Mostafa"
$ws.Cells.Item(117, 7).Value = "Dr. Amr Saeed, Dr. Enas Omran, Dr. Taqwa Mohammad"
$ws.Cells.Item(119, 7).Value = "Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Shorok Mohammad, Dr. Nahla, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Marina Sorial, Dr. Aya Hanafy"
$ws.Cells.Item(120, 7).Value = "Dr. Maryam Ashraf, Dr. Ola Abd Al-Fattah, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Gehad Salah, Dr. Yassmen Ahmad, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Aya Hanafy"
$ws.Cells.Item(121, 7).Value = "Dr. Eman Tantawi, Administrator, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(123, 7).Value = "D Wessam Atef, Dr. Amal Awwad, Dr. Sara Nabil, Dr. Nourhan Mohammad"
$ws.Cells.Item(124, 7).Value = "Dr. Amal Awwad, Dr. Sara Nabil, Dr. Nourhan Mohammad"
$ws.Cells.Item(128, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat, Dr. Mayar Ahmad Embaby"
$ws.Cells.Item(129, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef, Administrator"
$ws.Cells.Item(132, 7).Value = "Dr. Mariam Gamal Sanad, Dr. Sarah Mahdy"
$ws.Cells.Item(133, 7).Value = "Dr. Enas Omran, Dr. Afaf Abdallah, Dr. Nourham Mostafa"
$ws.Cells.Item(134, 7).Value = "Dr. Amr Saeed, Dr. Enas Omran, Dr. Taqwa Mohammad"
$ws.Cells.Item(136, 7).Value = "Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Shorok Mohammad, Dr. Nahla, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Marina Sorial, Dr. Aya Hanafy"
$ws.Cells.Item(137, 7).Value = "Dr. Maryam Ashraf, Dr. Ola Abd Al-Fattah, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Gehad Salah, Dr. Yassmen Ahmad, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Aya Hanafy"
$ws.Cells.Item(138, 7).Value = "Dr. Eman Tantawi, Administrator, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(142, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Basma Hamed, Dr. Arwa Al-Sayed"
$ws.Cells.Item(144, 7).Value = "Dr. Mohammad Safwat, Nourhan Mamdouh Hassan, Dr. Mayar Ahmad Embaby, Dr. Mariam Toma Gerges"
$ws.Cells.Item(145, 7).Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat, Dr. Mayar Ahmad Embaby"
$ws.Cells.Item(146, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat"
$ws.Cells.Item(148, 7).Value = "Dr. Mariam Gamal Sanad, Dr. Sarah Mahdy"
$ws.Cells.Item(150, 7).Value = "Dr. Youstina Ibrahim, Dr. Afaf Abdallah, Dr. Marian Samir"
$ws.Cells.Item(151, 7).Value = "Dr. Marian Samir, Dr. Rada Rabea, Administrator, Dr. Hana Amr, Dr. Nourhan Mohammad"
